$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

$row = 18

$ws.Cells.Item($row, 1).Value = "2025-08-29T06:32:52.639362"
$ws.Cells.Item($row, 2).Value = 5
$ws.Cells.Item($row, 3).Value = "全案件リスト"
$ws.Cells.Item($row, 4).Value = 80
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 5
